# EoCtUH (Efficiency of Conversion to Usable Heat) workbook update
#
# The "Efficiency" header on the EoCtUH sheet is clarified to make the
# dimensionless nature of the ratio explicit, matching the CHP efficiency
# values reported on the Data sheet.  The cell is widened (word-wrapped)
# and the row made taller so the longer label still reads cleanly.

$wb = $excel.ActiveWorkbook

$eoctuh = $wb.Worksheets.Item("EoCtUH")

$label = $eoctuh.Range("B1")
$label.Value = "Efficiency (dimensionless)"
$label.WrapText = $true
$eoctuh.Rows.Item(1).RowHeight = 45
$label.Select()

# Restore the originally active sheet/tab so the workbook's view state
# (which sheet is selected) is unchanged by this edit.
$wb.Worksheets.Item("About").Activate()
